$d = $word.ActiveDocument

# 1) "Protocol: Semantic Identifiers." -> "Protocol: Semantic Identifiers. URNs."
$old1 = "Protocol: Semantic Identifiers."
$new1 = "Protocol: Semantic Identifiers. URNs."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) "Protocol: Context Driven Interaction ... resolution)." -> new, longer description
$old2 = "Protocol: Context Driven Interaction REST P2P (SIDs CDI Dialogs: runat peer resolution addressable / browseable messages interactions embedded session semantics: events sourcing / history terms resolution)."
$new2 = "Protocol: Context Driven Interaction REST P2P (SIDs URNs: Resources  DCI Dialogs). Runat peer resolution addressable / browseable Messages interactions: request / response Message streams DCI dialogs. Embedded session semantics: event sourcing / history terms / roles resolution / navigation)."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3) "Messages: Case Classes. State flows." -> new text
$old3 = "Messages: Case Classes. State flows."
$new3 = "Messages: SIDs URNs Case Classes Statements. Statement Data Pattern Matching. State Flows: Reactive Events Messages."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
